# Adds two new columns, I0 (column I) and IF (column J), to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cells for the new columns, matching style of existing headers (row 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data values for rows 2..38 (I = I0 column, J = IF column)
$I0Values = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,7,1,5,1,5,1)
$IFValues = @(3,4,6,2,3,5,6,6,5,6,6,7,5,5,6,5,5,6,5,4,6,7,5,6,7,6,5,5,6,5,5,8,6,9,4,8,3)

for ($idx = 0; $idx -lt $I0Values.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $I0Values[$idx]
    $ws.Cells.Item($row, 10).Value = $IFValues[$idx]
}
